$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SanityTC")
$ws.Activate()

# "CaseToRun" column: these 4 sanity cases were re-used from another run,
# so they are no longer (re)run here -> flip from Y to N.
$ws.Range("B2").Value = "N"
$ws.Range("B3").Value = "N"
$ws.Range("B4").Value = "N"
$ws.Range("B5").Value = "N"

# "Pass/Fail/Skip" column: mark the reused/blocked cases as SKIP and give
# them a consistent status color (previously I2/I3 still carried the green
# "PASS" fill, which no longer matches the SKIP text).
$ws.Range("I2").Value = "SKIP"
$ws.Range("I3").Value = "SKIP"
$ws.Range("I4").Value = "SKIP"
$ws.Range("I5").Value = "SKIP"
$ws.Range("I2:I5").Interior.ColorIndex = 6

# "Failure Reason" column: note why the case results were skipped.
$ws.Range("J3").Value = " Sogo Account not present on page."
$ws.Range("J4").Value = " 05 not present on page."

# Restore the cursor/selection state left behind when the sheet was last saved.
$ws.Range("B13").Select()
